$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $NewValue)
    $escaped = $NewValue.Replace('"', '""')
    $rng = $ws.Range($CellRef)
    $rng.Formula = '="' + $escaped + '"'
    $rng.Copy()
    $rng.PasteSpecial(-4163)
}

Set-TextValue "D2" "265.22"
Set-TextValue "G2" "12"
Set-TextValue "D3" "22.86"
Set-TextValue "G3" "12"
Set-TextValue "D4" "6.251"
Set-TextValue "G4" "12"
Set-TextValue "G5" "12"
Set-TextValue "D6" "3.561"
Set-TextValue "G6" "12"
Set-TextValue "G7" "12"
Set-TextValue "D8" "1.358"
Set-TextValue "G8" "12"
Set-TextValue "D9" "0.8161"
Set-TextValue "G9" "12"
Set-TextValue "D10" "0.01353"
Set-TextValue "G10" "12"
Set-TextValue "D11" "0.1597"
Set-TextValue "G11" "12"
Set-TextValue "D12" "0.08191"
Set-TextValue "G12" "12"
Set-TextValue "G13" "12"
Set-TextValue "D14" "0.03139"
Set-TextValue "G14" "12"
Set-TextValue "D15" "0.09259"
Set-TextValue "G15" "12"
Set-TextValue "D16" "3.913"
Set-TextValue "G16" "12"
Set-TextValue "D17" "0.001692"
Set-TextValue "G17" "12"
Set-TextValue "D18" "0.04835"
Set-TextValue "G18" "12"
Set-TextValue "D19" "0.006235"
Set-TextValue "G19" "12"
Set-TextValue "B20" "BitKan"
Set-TextValue "C20" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D20" "0.001098"
Set-TextValue "E20" "19BitKanKAN"
Set-TextValue "G20" "12"
Set-TextValue "B21" "HotbitToken"
Set-TextValue "C21" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D21" "0.003229"
Set-TextValue "E21" "20HotbitTokenHTB"
Set-TextValue "G21" "12"
Set-TextValue "G22" "12"
Set-TextValue "D23" "3.696"
Set-TextValue "G23" "12"
Set-TextValue "D24" "2.261"
Set-TextValue "G24" "12"
Set-TextValue "D25" "0.3383"
Set-TextValue "G25" "12"
Set-TextValue "G26" "12"
Set-TextValue "D27" "0.0002682"
Set-TextValue "G27" "12"
Set-TextValue "G28" "12"
Set-TextValue "G29" "12"
Set-TextValue "G30" "12"
Set-TextValue "G31" "12"
Set-TextValue "G32" "12"
Set-TextValue "G33" "12"
Set-TextValue "G34" "12"
Set-TextValue "G35" "12"
Set-TextValue "G36" "12"
Set-TextValue "G37" "12"
Set-TextValue "G38" "12"
Set-TextValue "G39" "12"
Set-TextValue "D40" "0.04603"
Set-TextValue "G40" "12"
Set-TextValue "D41" "0.007220"
Set-TextValue "G41" "12"
Set-TextValue "D42" "0.1136"
Set-TextValue "G42" "12"
Set-TextValue "D43" "0.003400"
Set-TextValue "G43" "12"
Set-TextValue "G44" "12"
Set-TextValue "D45" "0.00006119"
Set-TextValue "G45" "12"
Set-TextValue "G46" "12"
Set-TextValue "D47" "0.7500"
Set-TextValue "G47" "12"
Set-TextValue "D48" "0.1973"
Set-TextValue "G48" "12"
Set-TextValue "G49" "12"
Set-TextValue "G50" "12"
Set-TextValue "G51" "12"

$excel.CutCopyMode = 0

